$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-07 21:18:40"
$ws.Range("E3").Value = "2026-02-07 21:18:42"
$ws.Range("E4").Value = "2026-02-07 21:18:45"
$ws.Range("J4").Value = "1004.0 hPa"
$ws.Range("N4").Value = "9.3 °C 20:59 TU"
$ws.Range("E5").Value = "2026-02-07 21:18:47"
$ws.Range("L5").Value = "47.5 km/h - 64º 20:39 TU"
$ws.Range("E6").Value = "2026-02-07 21:18:50"
$ws.Range("J6").Value = "1003.9 hPa"
$ws.Range("E7").Value = "2026-02-07 21:18:53"
$ws.Range("E8").Value = "2026-02-07 21:18:56"
$ws.Range("E9").Value = "2026-02-07 21:18:58"
$ws.Range("N9").Value = "4.3 °C 20:54 TU"
$ws.Range("O9").Value = "11.0 °C"
$ws.Range("E10").Value = "2026-02-07 21:19:01"
$ws.Range("O10").Value = "8.4 °C"
$ws.Range("E11").Value = "2026-02-07 21:19:04"
$ws.Range("H11").Value = "'86%"
$ws.Range("O11").Value = "3.5 °C"
$ws.Range("E12").Value = "2026-02-07 21:19:06"
$ws.Range("H12").Value = "'86%"
$ws.Range("N12").Value = "5.5 °C 20:43 TU"
$ws.Range("O12").Value = "10.3 °C"
$ws.Range("E13").Value = "2026-02-07 21:19:09"
$ws.Range("E14").Value = "2026-02-07 21:19:12"
$ws.Range("O14").Value = "11.9 °C"
$ws.Range("E15").Value = "2026-02-07 21:19:14"
$ws.Range("N15").Value = "4.6 °C 20:59 TU"
$ws.Range("O15").Value = "10.6 °C"
$ws.Range("E16").Value = "2026-02-07 21:19:17"
$ws.Range("H16").Value = "'59%"
$ws.Range("E17").Value = "2026-02-07 21:19:20"
$ws.Range("L17").Value = "65.2 km/h - 242º 20:59 TU"
$ws.Range("E18").Value = "2026-02-07 21:19:22"
$ws.Range("E19").Value = "2026-02-07 21:19:25"
$ws.Range("E20").Value = "2026-02-07 21:19:28"
$ws.Range("O20").Value = "-5.7 °C"
$ws.Range("E21").Value = "2026-02-07 21:19:30"
$ws.Range("H21").Value = "'82%"
$ws.Range("J21").Value = "1006.1 hPa"
$ws.Range("E22").Value = "2026-02-07 21:19:33"
$ws.Range("L22").Value = "24.5 km/h - 203º 20:45 TU"
$ws.Range("O22").Value = "-6.2 °C"
$ws.Range("E23").Value = "2026-02-07 21:19:36"
$ws.Range("L23").Value = "39.2 km/h - 145º 20:59 TU"
$ws.Range("E24").Value = "2026-02-07 21:19:38"
$ws.Range("I24").Value = "0.5 mm"
$ws.Range("O24").Value = "7.6 °C"
$ws.Range("E25").Value = "2026-02-07 21:19:41"
$ws.Range("O25").Value = "-3.9 °C"
$ws.Range("E26").Value = "2026-02-07 21:19:43"
$ws.Range("E27").Value = "2026-02-07 21:19:46"
$ws.Range("E28").Value = "2026-02-07 21:19:48"
$ws.Range("O28").Value = "8.6 °C"
$ws.Range("E29").Value = "2026-02-07 21:19:51"
$ws.Range("H29").Value = "'72%"
$ws.Range("E30").Value = "2026-02-07 21:19:53"
$ws.Range("H30").Value = "'80%"
$ws.Range("N30").Value = "5.2 °C 20:55 TU"
$ws.Range("O30").Value = "9.9 °C"
$ws.Range("E31").Value = "2026-02-07 21:19:56"
$ws.Range("E32").Value = "2026-02-07 21:19:58"
$ws.Range("H32").Value = "'77%"
$ws.Range("E33").Value = "2026-02-07 21:20:01"
$ws.Range("E34").Value = "2026-02-07 21:20:04"
$ws.Range("K34").Value = "10.5 MJ/m2"
$ws.Range("O34").Value = "-2.4 °C"
$ws.Range("E35").Value = "2026-02-07 21:20:06"
$ws.Range("E36").Value = "2026-02-07 21:20:09"
$ws.Range("E37").Value = "2026-02-07 21:20:12"
$ws.Range("J37").Value = "1005.9 hPa"
$ws.Range("E38").Value = "2026-02-07 21:20:14"
$ws.Range("O38").Value = "12.1 °C"
$ws.Range("E39").Value = "2026-02-07 21:20:17"
$ws.Range("H39").Value = "'67%"
$ws.Range("E40").Value = "2026-02-07 21:20:20"
$ws.Range("E41").Value = "2026-02-07 21:20:22"
$ws.Range("E42").Value = "2026-02-07 21:20:25"
$ws.Range("H42").Value = "'75%"
$ws.Range("E43").Value = "2026-02-07 21:20:27"
$ws.Range("O43").Value = "7.7 °C"
$ws.Range("E44").Value = "2026-02-07 21:20:30"
$ws.Range("O44").Value = "-4.4 °C"
$ws.Range("E45").Value = "2026-02-07 21:20:33"
$ws.Range("H45").Value = "'61%"
$ws.Range("J45").Value = "1004.7 hPa"
$ws.Range("O45").Value = "4.1 °C"
$ws.Range("E46").Value = "2026-02-07 21:20:35"
$ws.Range("J46").Value = "1007.6 hPa"
